$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 already exists (previously blank, formatted) - copy its formatting (A:J, not K)
# onto the new row 5 so the new row matches the sheet's existing look-and-feel.
$ws.Range("A6:J6").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate new row 5 - "Crumpet" facility
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Populate previously-blank row 6 - "Scone" facility
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Match the author's final selection/view: range A5:J6 selected, active cell A5.
$ws.Activate() | Out-Null
$ws.Range("A5:J6").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
